$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.270.86"
$ws.Range("E2").Value = "  -2.15%  "

$ws.Range("D3").Value = "2.238.13"
$ws.Range("E3").Value = "  -1.97%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "'230.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.74%  "

$ws.Range("E6").Value = "  -0.10%  "

$ws.Range("D7").Value = "'63.59"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.10%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("D9").Value = "'0.440"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.46%  "

$ws.Range("D10").Value = "'0.0951"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -8.90%  "

$ws.Range("D11").Value = "'56.47"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.56%  "

$ws.Range("D12").Value = "'27.45"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.56%  "

$ws.Range("E13").Value = "  -1.58%  "

$ws.Range("D14").Value = "2.569.37"
$ws.Range("E14").Value = "  -2.12%  "

$ws.Range("D15").Value = "'15.41"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.32%  "

$ws.Range("D16").Value = "'6.04"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.19%  "

$ws.Range("D17").Value = "'0.824"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.19%  "

$ws.Range("D18").Value = "2.234.10"
$ws.Range("E18").Value = "  -2.18%  "

$ws.Range("D19").Value = "43.173.00"
$ws.Range("E19").Value = "  -2.08%  "

$ws.Range("D20").Value = "0.0₃0963"
$ws.Range("E20").Value = "  -6.44%  "

$ws.Range("E21").Value = "  -1.85%  "

$ws.Range("E22").Value = "  -0.33%  "

$ws.Range("D23").Value = "'245.61"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.52%  "

$ws.Range("D24").Value = "'0.999"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.06%  "

$ws.Range("D25").Value = "'3.68"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +29.97%  "

$ws.Range("D26").Value = "'2.41"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.11%  "

$ws.Range("E27").Value = "  -4.48%  "

$ws.Range("D28").Value = "'9.71"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.04%  "

$ws.Range("D29").Value = "'173.23"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.84%  "

$ws.Range("D30").Value = "'21.46"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.73%  "

$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").Value = "'1.40"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.44%  "

$ws.Range("B32").Value = "Kaspa"
$ws.Range("C32").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D32").Value = "'0.128"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -7.78%  "

$ws.Range("E33").Value = "  -0.33%  "

$ws.Range("D34").Value = "'4.94"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.83%  "

$ws.Range("E35").Value = "  -2.25%  "

$ws.Range("E36").Value = "  -2.10%  "

$ws.Range("D37").Value = "'3.58"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -7.63%  "

$ws.Range("E38").Value = "  -8.39%  "

$ws.Range("E39").Value = "  -3.95%  "

$ws.Range("D40").Value = "'0.0249"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.04%  "

$ws.Range("D41").Value = "'1.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.13%  "

$ws.Range("E42").Value = "  +1.09%  "

$ws.Range("E43").Value = "  -0.15%  "

$ws.Range("D44").Value = "'16.92"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.67%  "

$ws.Range("D45").Value = "'96.39"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.91%  "

$ws.Range("D46").Value = "'0.0941"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.36%  "

$ws.Range("D47").Value = "'0.000208"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.13%  "

$ws.Range("E48").Value = "  -2.67%  "

$ws.Range("D49").Value = "1.441.85"
$ws.Range("E49").Value = "  -2.16%  "

$ws.Range("D50").Value = "'9.89"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.20%  "

$ws.Range("D51").Value = "'2.26"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.98%  "
